# Update "想去人数" (F column) figures on the "展览" and "全部类型" sheets
# to reflect the latest generated output (gh-pages rebuild at 456a3b4).

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 11665
    3  = 11249
    6  = 1019
    11 = 10719
    12 = 4145
    20 = 443
    21 = 11129
    22 = 10908
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
